# Auto-generated edit script: apply per-cell value updates described by the diff
# against the 'Sheets/Aegis_Profits.xlsx' workbook (8 worksheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 11182.5
$ws.Range("I51").Value = 16343
$ws.Range("J51").Value = 3957.8
$ws.Range("K51").Value = 16343
$ws.Range("L51").Value = 3957.8
$ws.Range("M51").Value = -15859
$ws.Range("N51").Value = -4925.8

$ws.Range("H55").Value = 569.3
$ws.Range("I55").Value = 540
$ws.Range("J55").Value = 613.25
$ws.Range("K55").Value = 540
$ws.Range("L55").Value = 613.25
$ws.Range("M55").Value = -326
$ws.Range("N55").Value = -1041.25

$ws.Range("H64").Value = 69693.2
$ws.Range("I64").Value = 202400
$ws.Range("J64").Value = 3339.8
$ws.Range("K64").Value = 202400
$ws.Range("L64").Value = 3339.8
$ws.Range("M64").Value = -202152
$ws.Range("N64").Value = -3835.8

$ws.Range("H67").Value = 69693.2
$ws.Range("I67").Value = 202400
$ws.Range("J67").Value = 3339.8
$ws.Range("K67").Value = 202400
$ws.Range("L67").Value = 3339.8
$ws.Range("M67").Value = -201542
$ws.Range("N67").Value = -5055.8

$ws.Range("H137").Value = 1297.5333
$ws.Range("I137").Value = 1097.4286
$ws.Range("J137").Value = 1997.9
$ws.Range("K137").Value = 3292.2858
$ws.Range("L137").Value = 5993.700000000001
$ws.Range("M137").Value = -742.2857999999997
$ws.Range("N137").Value = -11093.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20700.404
$ws.Range("I32").Value = 3632.3247
$ws.Range("K32").Value = 3632.3247
$ws.Range("M32").Value = -3345.3247

$ws.Range("H74").Value = 1932.4117
$ws.Range("I74").Value = 786.4167
$ws.Range("J74").Value = 4682.8
$ws.Range("K74").Value = 786.4167
$ws.Range("L74").Value = 4682.8
$ws.Range("M74").Value = 87.58330000000001
$ws.Range("N74").Value = -6430.8

$ws.Range("H77").Value = 1932.4117
$ws.Range("I77").Value = 786.4167
$ws.Range("J77").Value = 4682.8
$ws.Range("K77").Value = 3932.0835
$ws.Range("L77").Value = 23414
$ws.Range("M77").Value = 435.9165000000003
$ws.Range("N77").Value = -32150

$ws.Range("H113").Value = 33800
$ws.Range("J113").Value = 33800
$ws.Range("L113").Value = 33800
$ws.Range("N113").Value = -42478

$ws.Range("H132").Value = 4005.5
$ws.Range("I132").Value = 4168.9165
$ws.Range("K132").Value = 12506.7495
$ws.Range("M132").Value = -9976.749500000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 862.375
$ws.Range("I94").Value = 651.75
$ws.Range("K94").Value = 651.75
$ws.Range("M94").Value = -200.75

$ws.Range("H99").Value = 1485.7368
$ws.Range("I99").Value = 1357.2727
$ws.Range("J99").Value = 1662.375
$ws.Range("K99").Value = 1357.2727
$ws.Range("L99").Value = 1662.375
$ws.Range("M99").Value = 140.7273
$ws.Range("N99").Value = -4658.375

$ws.Range("H134").Value = 3132.724
$ws.Range("I134").Value = 3144.6072
$ws.Range("K134").Value = 9433.821599999999
$ws.Range("M134").Value = -6898.821599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H134").Value = 1241.6842
$ws.Range("I134").Value = 1092.2858
$ws.Range("J134").Value = 1660
$ws.Range("K134").Value = 3276.8574
$ws.Range("L134").Value = 4980
$ws.Range("M134").Value = -741.8574000000003
$ws.Range("N134").Value = -10050

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 15342.286
$ws.Range("I80").Value = 400
$ws.Range("J80").Value = 17832.666
$ws.Range("K80").Value = 1200
$ws.Range("L80").Value = 53497.99800000001
$ws.Range("M80").Value = -264
$ws.Range("N80").Value = -55369.99800000001

$ws.Range("H83").Value = 15342.286
$ws.Range("I83").Value = 400
$ws.Range("J83").Value = 17832.666
$ws.Range("K83").Value = 3600
$ws.Range("L83").Value = 160493.994
$ws.Range("M83").Value = 1080
$ws.Range("N83").Value = -169853.994

$ws.Range("H113").Value = 744.5
$ws.Range("J113").Value = 737.1111
$ws.Range("L113").Value = 2211.3333
$ws.Range("N113").Value = -6551.3333

$ws.Range("H131").Value = 819.59
$ws.Range("I131").Value = 447
$ws.Range("J131").Value = 860.9888999999999
$ws.Range("K131").Value = 1341
$ws.Range("L131").Value = 2582.9667
$ws.Range("M131").Value = 3699
$ws.Range("N131").Value = -12662.9667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 5237000
$ws.Range("J11").Value = 75999
$ws.Range("L11").Value = 75999
$ws.Range("N11").Value = -76277

$ws.Range("H70").Value = 58690.434
$ws.Range("I70").Value = 121879.414
$ws.Range("J70").Value = 4979.8
$ws.Range("K70").Value = 121879.414
$ws.Range("L70").Value = 4979.8
$ws.Range("M70").Value = -121609.414
$ws.Range("N70").Value = -5519.8

$ws.Range("H73").Value = 58690.434
$ws.Range("I73").Value = 121879.414
$ws.Range("J73").Value = 4979.8
$ws.Range("K73").Value = 121879.414
$ws.Range("L73").Value = 4979.8
$ws.Range("M73").Value = -120943.414
$ws.Range("N73").Value = -6851.8

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H126").Value = 3274.2
$ws.Range("I126").Value = 3082.4443
$ws.Range("K126").Value = 9247.332900000001
$ws.Range("M126").Value = -6777.332900000001

$ws.Range("H132").Value = 1855.3103
$ws.Range("I132").Value = 1363.3684
$ws.Range("J132").Value = 2790
$ws.Range("K132").Value = 4090.1052
$ws.Range("L132").Value = 8370
$ws.Range("M132").Value = -1560.1052
$ws.Range("N132").Value = -13430

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 40674.08
$ws.Range("I16").Value = 59302.94
$ws.Range("J16").Value = 1087.75
$ws.Range("K16").Value = 59302.94
$ws.Range("L16").Value = 1087.75
$ws.Range("M16").Value = -59132.94
$ws.Range("N16").Value = -1427.75

$ws.Range("H22").Value = 821.2778
$ws.Range("J22").Value = 727.4286
$ws.Range("L22").Value = 727.4286
$ws.Range("N22").Value = -1317.4286

$ws.Range("H27").Value = 821.2778
$ws.Range("J27").Value = 727.4286
$ws.Range("L27").Value = 727.4286
$ws.Range("N27").Value = -941.4286

$ws.Range("H46").Value = 2949.7778
$ws.Range("J46").Value = 2949.7778
$ws.Range("L46").Value = 2949.7778
$ws.Range("N46").Value = -3325.7778

$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H93").Value = 1997
$ws.Range("I93").Value = 2396
$ws.Range("J93").Value = 999.5
$ws.Range("K93").Value = 2396
$ws.Range("L93").Value = 999.5
$ws.Range("M93").Value = -1148
$ws.Range("N93").Value = -3495.5

$ws.Range("H132").Value = 3981.682
$ws.Range("I132").Value = 3981.682
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11945.046
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9415.045999999998
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H81").Value = 401091.8
$ws.Range("I81").Value = 501000
$ws.Range("J81").Value = 334486.34
$ws.Range("K81").Value = 1002000
$ws.Range("L81").Value = 668972.6800000001
$ws.Range("M81").Value = -1000939
$ws.Range("N81").Value = -671094.6800000001

$ws.Range("H84").Value = 401091.8
$ws.Range("I84").Value = 501000
$ws.Range("J84").Value = 334486.34
$ws.Range("K84").Value = 5010000
$ws.Range("L84").Value = 3344863.4
$ws.Range("M84").Value = -5004696
